# Generate Report for Archive
# - Status text "Ready for handoff" -> "In Translation" on all three sheets
#   (Overview!E2:F2, zh-cn!C2, de-de!C2)
# - Narrow the Status column(s) to match the shorter text
#   (Overview cols E:F, zh-cn col C, de-de col C)

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Update the status values
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus

# Resize the columns that held the status text to fit the new, shorter value
$wsOverview.Columns("E:F").ColumnWidth = 12.5
$wsZhCn.Columns("C:C").ColumnWidth = 12.5
$wsDeDe.Columns("C:C").ColumnWidth = 12.5
